{"js": "// Apply the CyberCatz copy refresh described by the commit diff:\n//  - retitle the review (heading + the bold recap line near the end)\n//  - rewrite the four \"What we like\" bullets\n//  - drop the \"Wild symbol does not replace...\" bullet from \"What we don't like\"\n//  - rewrite the italic summary/meta line at the end\n\nconst body = context.document.body;\n\n// Exact whole-paragraph text replacements: oldText -> newText.\n// Using matchCase + whole-paragraph search keeps this robust regardless of\n// how the run(s) inside each paragraph are split.\nconst replacements = [\n  [\n    \"Play CyberCatz Free: Intergalactic Slot Game Review\",\n    \"Play CyberCatz Online Slot Game for Free\",\n  ],\n  [\n    \"Unique and immersive graphics and soundtrack\",\n    \"Immersive intergalactic theme\",\n  ],\n  [\n    \"Free spins feature with up to 30 spins\",\n    \"Wild symbol for maximizing winnings\",\n  ],\n  [\n    \"Strategic use of Wild symbol\",\n    \"Free spins feature for additional chances of winning\",\n  ],\n  [\n    \"Exciting Cyber City Bonus game\",\n    \"Cyber City Bonus game for extra prizes\",\n  ],\n  [\n    \"Discover the features of CyberCatz in this slot game review. Play for free and enjoy the unique graphics and free spins feature for a chance to win big.\",\n    \"Read our review of the CyberCatz slot game and play for free to win big.\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// Remove the \"Wild symbol does not replace bonus icons or Free Spin\n// symbols\" bullet entirely (the whole paragraph is dropped).\nconst removed = body.search(\n  \"Wild symbol does not replace bonus icons or Free Spin symbols\",\n  { matchCase: true }\n);\nremoved.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < removed.items.length; i++) {\n  removed.items[i].paragraphs.getFirst().delete();\n}\nawait context.sync();\n", "ps1": "# Apply the CyberCatz copy refresh described by the commit diff:\n#  - retitle the review (heading + the bold recap line near the end)\n#  - rewrite the four \"What we like\" bullets\n#  - drop the \"Wild symbol does not replace...\" bullet from \"What we don't like\"\n#  - rewrite the italic summary/meta line at the end\n\n$d = $word.ActiveDocument\n\n# wdReplaceAll = 2 ; MatchCase = $true keeps these surgical, whole-phrase swaps.\nfunction Replace-AllText($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.Execute(\n        $findText,      # FindText\n        $true,          # MatchCase\n        $true,          # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $replaceText,   # ReplaceWith\n        2               # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n\n# Title: occurs twice (Heading1 at the top + the bold recap paragraph near\n# the end) \u2014 Replace All handles both hits in a single call.\nReplace-AllText \"Play CyberCatz Free: Intergalactic Slot Game Review\" \"Play CyberCatz Online Slot Game for Free\"\n\n# \"What we like\" bullets.\nReplace-AllText \"Unique and immersive graphics and soundtrack\" \"Immersive intergalactic theme\"\nReplace-AllText \"Free spins feature with up to 30 spins\" \"Wild symbol for maximizing winnings\"\nReplace-AllText \"Strategic use of Wild symbol\" \"Free spins feature for additional chances of winning\"\nReplace-AllText \"Exciting Cyber City Bonus game\" \"Cyber City Bonus game for extra prizes\"\n\n# Italic summary line at the very end of the document.\nReplace-AllText \"Discover the features of CyberCatz in this slot game review. Play for free and enjoy the unique graphics and free spins feature for a chance to win big.\" \"Read our review of the CyberCatz slot game and play for free to win big.\"\n\n# Drop the \"Wild symbol does not replace bonus icons or Free Spin symbols\"\n# bullet from \"What we don't like\" entirely (whole paragraph removed).\n$target = \"Wild symbol does not replace bonus icons or Free Spin symbols\"\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -ceq $target) {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
